$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("C1").Value = "CurrentSeason"
$ws.Range("D1").Value = "CurrentStageName"
$ws.Range("E1").Value = "BudBurstDOY"
$ws.Range("F1").Value = "FloweringDOY"
$ws.Range("G1").Value = "VeraisonDOY"
